$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Novischförman" column (D) is removed entirely; the former column E
# ("PQE") shifts left into D, and all following columns shift left too.
$ws.Columns("D").Delete()

# Row 12 was a stray leftover row (single cell with the value "s"); it is
# removed, shrinking the used range back down to A1:D8.
$ws.Rows(12).Delete()

# Move the active selection to the new end-of-data cell.
$ws.Range("E12").Select() | Out-Null

# Touch page setup (orientation) so the sheet gets a <pageSetup> entry.
$ps = $ws.PageSetup
$ps.Orientation = 1
